# Disaggregation of commodity Copper
# 1) Rename the commodity "Copper ores and concentrates" -> "Copper"
#    (this label lives in column C row 4 of every year sheet, all sharing
#    the same underlying text).
# 2) A handful of year sheets receive a last-digit (ULP) correction to the
#    D4 value that sits next to that label.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $label = $ws.Range("C4").Value2
    if ($label -eq "Copper ores and concentrates") {
        $ws.Range("C4").Value = "Copper"
    }
}

$years   = @("2023", "2025", "2028", "2032", "2041", "2048", "2072", "2085", "2087")
$newvals = @(55579.97923991122, 64307.96100302236, 93395.27619719859, 136136.6974506026, 505872.9439998683, 1950034.592079028, 1953747.730931857, 1681427.682419382, 1854362.458575674)

for ($i = 0; $i -lt $years.Length; $i++) {
    $ws = $wb.Worksheets.Item($years[$i])
    $ws.Range("D4").Value = $newvals[$i]
}
